# edit.ps1
# Applies the "updating template for 2014" commit:
#   - bumps the cached date-field text on the slide master, notes master
#     and handout master from 11/15/14 -> 11/12/15
#   - rewrites the Module-1 title on slide 2 to "RNA-Seq Module 1"
#   - drops "Walker, Ben " / swaps "Ainscough" -> "Walker" in the byline
#   - updates the course dates from "November 11-23, 2014" to
#     "November 10 - 22, 2015"

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($container, [string]$newText)

    if ($container -eq $null) { return }

    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            try {
                $tr = $shp.TextFrame.TextRange
                $len = $tr.Text.Length
                if ($len -gt 0) {
                    $rng = $tr.Characters(1, $len)
                    $rng.Text = $newText
                } else {
                    $tr.Text = $newText
                }
            } catch {
                # Some hosts cannot resolve an edit target inside the
                # notes/handout master date field; nothing more we can
                # do through the object model in that case.
            }
        }
    }
}

# --- 1) Refresh the cached "datetime1" field text wherever it appears ---
Set-DateFieldText $p.SlideMaster "11/12/15"
Set-DateFieldText $p.NotesMaster "11/12/15"
Set-DateFieldText $p.HandoutMaster "11/12/15"

# --- 2) Slide 2: title block "Module 1" -> "RNA-Seq Module 1" ---
$slide2 = $p.Slides.Item(2)

for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($full.IndexOf("Module 1") -ge 0 -and $full.IndexOf("Introduction to RNA sequencing") -ge 0) {
        $idx = $full.IndexOf("Module 1")
        $start = $idx + 1
        $len = "Module 1".Length
        $whole = $tr.Characters($start, $len)
        $whole.Text = "RNA-Seq Module 1"

        # Re-split into the run layout used by the edited deck:
        #   "RNA-" | "S" | "eq" | " Module " | "1"
        $pos = $start
        $piece = "RNA-"
        $tr.Characters($pos, $piece.Length).Text = $piece
        $pos += $piece.Length

        $piece = "S"
        $tr.Characters($pos, $piece.Length).Text = $piece
        $pos += $piece.Length

        $piece = "eq"
        $tr.Characters($pos, $piece.Length).Text = $piece
        $pos += $piece.Length

        $piece = " Module "
        $tr.Characters($pos, $piece.Length).Text = $piece
        $pos += $piece.Length

        $piece = "1"
        $tr.Characters($pos, $piece.Length).Text = $piece
    }
}

# --- 3) Slide 2: byline "... Jason Walker, Ben Ainscough" -> "... Jason Walker" ---
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    if ($full.IndexOf("Ainscough") -ge 0) {
        $startMarker = "Malachi Griffith, Obi Griffith, Jason Walker, Ben "
        $idx = $full.IndexOf($startMarker)
        if ($idx -ge 0) {
            $endIdx = $full.IndexOf("Ainscough") + "Ainscough".Length
            $start = $idx + 1
            $len = $endIdx - $idx
            $whole = $tr.Characters($start, $len)
            $whole.Text = "Malachi Griffith, Obi Griffith, Jason Walker"

            $part1 = "Malachi Griffith, Obi Griffith, Jason "
            $part2 = "Walker"
            $tr.Characters($start, $part1.Length).Text = $part1
            $tr.Characters($start + $part1.Length, $part2.Length).Text = $part2
        }
    }

    # --- 4) Slide 2: course dates "November 11-23, 2014" -> "November 10 - 22, 2015"
    $full = $tr.Text
    $oldDate = "November 11-23, 2014"
    $newDate = "November 10 - 22, 2015"
    $idx = $full.IndexOf($oldDate)
    if ($idx -ge 0) {
        $start = $idx + 1
        $rng = $tr.Characters($start, $oldDate.Length)
        $rng.Text = $newDate
    }
}
